# Rebuild the single "questions = [...]" payload as pretty-printed JSON
# (double quotes, 4-space indent, \uXXXX-escaped smart quotes/apostrophes)
# instead of the old Python-repr-with-single-quotes string.
$newText = @'
questions = [
    {
        "title": "You work in the technology department of a department store. You are conducting a needs assessment with a client to determine the best camera to show them. The client mentions that they will be using the product on family trips.  How should you respond?",
        "ques_type": 2,
        "options": [
            "\"What sort of trips do you take with your family?\"",
            "\"This would be great for you because it has a long battery life and is difficult to break.\"",
            "\"These are the five cameras that are best for family trips.\"",
            " \u201cThis camera here is my favorite, as it has a high-quality zoom.\""
        ],
        "score": "\"What sort of trips do you take with your family?\""
    },
    {
        "title": "You work as a sales consultant for ABC computers. While talking to the customer about one of your laptops, the customer mentions a competing laptop from DEF, your direct competitor, that they feel has some advantages over yours.  What should you say?",
        "ques_type": 2,
        "options": [
            "\u201cI understand you\u2019re looking for the best laptop. The ABC model has a faster processing speed.\u201d",
            "\u201cI understand you\u2019re looking for the best laptop. DEF makes a great laptop as well. Either product would be a great choice.\u201d",
            "\u201cI understand you\u2019re looking for the best laptop. Let\u2019s look at a side-by-side comparison of both products so you can make an informed decision.\u201d",
            "\u201cI understand you\u2019re looking for the best laptop. Let me tell you all the ways I feel ABC\u2019s laptop is superior to DEF\u2019s.\u201d"
        ],
        "score": "\u201cI understand you\u2019re looking for the best laptop. Let\u2019s look at a side-by-side comparison of both products so you can make an informed decision.\u201d"
    },
    {
        "title": "You are working at a car dealership where research done by management shows a low rate of return business. To increase your repeat and referral business, you have decided to create a Facebook business page for your customers.  Which of the following actions should most effectively help you achieve your goal?",
        "ques_type": 2,
        "options": [
            "Tag them in your feed about sales specials on accessories they might like for their vehicle.",
            "Invite them to a group page where you share current news about the company.",
            "Like their most popular posts.",
            "Send them direct, personalized messages on important dates, such as their birthday and wedding anniversary."
        ],
        "score": "Send them direct, personalized messages on important dates, such as their birthday and wedding anniversary."
    },
    {
        "title": "True or false: When speaking to potential customers, you should always try to speak the most so you can control the conversation.",
        "ques_type": 11,
        "options": [
            "true",
            "false"
        ],
        "score": "False"
    }
]
'@

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old sheet had two rows: A1 (a bold/bordered/centered "0" placeholder)
# and A2 (the shared-string JSON blob). Drop row 1 entirely so the JSON
# blob becomes A1, with the plain default style (no bold/border/alignment).
$ws.Rows.Item(1).Delete()
$ws.Range("A1").Value = $newText

# Make sure no custom row-height sticks around from the long single-line
# text (Excel autosizes rows on entry of long strings).
$ws.Rows.Item(1).AutoFit()
